$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 70
$ws.Range("B3").Value = 288.7
$ws.Range("C3").Value = 195.1
$ws.Range("C4").Value = 284.5
$ws.Range("C11").Value = 259.9
$ws.Range("C13").Value = 336.3
$ws.Range("C17").Value = 44.8
$ws.Range("C19").Value = -106.5
$ws.Range("C23").Value = 21.5
